$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.962.68'
$ws.Range("E2").Value = '  -1.79%  '
$ws.Range("D3").Value = '3.383.44'
$ws.Range("E3").Value = '  -0.79%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '574.26'
$ws.Range("E5").Value = '  -1.40%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '136.63'
$ws.Range("E6").Value = '  -1.81%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").Value = '3.383.25'
$ws.Range("E8").Value = '  -0.79%  '
$ws.Range("E9").Value = '  -1.48%  '
$ws.Range("E10").Value = '  +2.11%  '
$ws.Range("E11").Value = '  -4.05%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.382'
$ws.Range("E12").Value = '  -2.95%  '
$ws.Range("D13").Value = '3.957.77'
$ws.Range("E13").Value = '  -0.86%  '
$ws.Range("E14").Value = '  +0.83%  '
$ws.Range("B15").Value = 'ShibaInu'
$ws.Range("C15").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000172'
$ws.Range("E15").Value = '  -4.10%  '
$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").Value = '3.377.98'
$ws.Range("E16").Value = '  -1.62%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '25.40'
$ws.Range("E17").Value = '  -0.80%  '
$ws.Range("D18").Value = '61.126.30'
$ws.Range("E18").Value = '  -1.67%  '
$ws.Range("E19").Value = '  -2.81%  '
$ws.Range("E20").Value = '  -1.71%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.35'
$ws.Range("E21").Value = '  -2.39%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '376.22'
$ws.Range("E22").Value = '  -5.21%  '
$ws.Range("D23").Value = '3.517.91'
$ws.Range("E23").Value = '  -0.96%  '
$ws.Range("E24").Value = '  -3.12%  '
$ws.Range("E25").Value = '  +0.36%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000126'
$ws.Range("E26").Value = '  -4.54%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '70.96'
$ws.Range("E27").Value = '  -1.04%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.182'
$ws.Range("E28").Value = '  +12.15%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.65'
$ws.Range("E29").Value = '  -1.52%  '
$ws.Range("E30").Value = '  +0.01%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.40'
$ws.Range("E31").Value = '  -4.31%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.07'
$ws.Range("E32").Value = '  -2.52%  '
$ws.Range("E33").Value = '  -2.14%  '
$ws.Range("E34").Value = '  -0.03%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '23.50'
$ws.Range("E35").Value = '  -0.51%  '
$ws.Range("E36").Value = '  -4.71%  '
$ws.Range("E37").Value = '  -3.49%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.80'
$ws.Range("E38").Value = '  -2.35%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '164.29'
$ws.Range("E39").Value = '  -0.36%  '
$ws.Range("E40").Value = '  -4.88%  '
$ws.Range("E41").Value = '  -0.07%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '25.28'
$ws.Range("E42").Value = '  +0.66%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.774'
$ws.Range("E43").Value = '  -2.04%  '
$ws.Range("E44").Value = '  -0.11%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.70'
$ws.Range("E45").Value = '  -5.88%  '
$ws.Range("E46").Value = '  -6.77%  '
$ws.Range("E47").Value = '  -2.81%  '
$ws.Range("D48").Value = '2.487.97'
$ws.Range("E48").Value = '  +5.39%  '
$ws.Range("E49").Value = '  -2.25%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '22.85'
$ws.Range("E50").Value = '  -2.87%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.44'
$ws.Range("E51").Value = '  +3.55%  '
